$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows after the existing data row 17 (shifts old rows 18-23 down to 20-25)
$ws.Rows("18:19").Insert()

# Row 19 should inherit the "last row" (thicker bottom border) formatting that
# row 17 currently has, since it becomes the new last row of the table.
$ws.Range("B17:J17").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)

# Row 18 (new) and row 17 (existing, no longer last) both use the "middle row" formatting from row 16.
$ws.Range("B16:J16").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Fill in row 18: ZAMIR GARCIA CARMONA, new period 2508
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1143355980"
$ws.Range("D18").Value = "ZAMIR GARCIA CARMONA"
$ws.Range("E18").Value = "2508"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500

# Fill in row 19: FRAY DE JESUS ESCORCIA SALAS, new period 2508
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1049533731"
$ws.Range("D19").Value = "FRAY DE JESUS ESCORCIA SALAS"
$ws.Range("E19").Value = "2508"
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500

# Update summary totals: Valor Mora doubled, Cant. Periodos now 2
$ws.Range("E11").Value = 227760
$ws.Range("F13").Value = 2
